# Applies the edits described in the commit: rewrites five single-message /
# paired-message phishing bodies, and swaps the order of the "Felicitations
# Sue" / "FoodFlow" paired-message-B content blocks.

$d = $word.ActiveDocument

# Single message 1: Gary Leal fitness tracker -> Mr. Jeter / Global Activist Network
$d.Paragraphs.Item(8).Range.Text = "Dear Mr. Jeter," + [char]11 + "" + [char]11 + " We are reaching out to you from the Customer Care team at Global Activist Network - a platform that connects passionate individuals like you to the causes they care about. We recently conducted a system-wide audit and it appears that there may have been an error with your last donation transaction." + [char]11 + "" + [char]11 + " To rectify the situation, we need to verify your credit card details. Please reply to this text with your full name as it appears on the card, the 16-digit card number, expiry date, and CVV. Your prompt response would help us ensure that your contributions continue to reach those who need them most, without interruption. " + [char]11 + "" + [char]11 + " We apologize for any inconvenience this may cause and appreciate your understanding and support. " + [char]11 + "" + [char]11 + " Best," + [char]11 + " Customer Care Team" + [char]11 + " Global Activist Network`""

# Single message 2: Mr. Vega credit card renewal -> Justin / Nature Magazine subscription
$d.Paragraphs.Item(11).Range.Text = "Subject: Justin, Your Subscription is Expiring Soon" + [char]11 + "" + [char]11 + " Dear Justin," + [char]11 + "" + [char]11 + " We hope this message finds you well. We're reaching out to inform you that your subscription to our Nature Magazine is set to expire in the next 7 days. " + [char]11 + "" + [char]11 + " To ensure you continue receiving the latest issues without interruption, we kindly request that you update your payment information by clicking the button below and providing your credit card details." + [char]11 + "" + [char]11 + " This is a quick and easy process that will allow us to seamlessly renew your subscription. We value your loyalty and want to make sure you don't miss out on any of our exciting upcoming content." + [char]11 + "" + [char]11 + " Please click here to update your payment details:" + [char]11 + "" + [char]11 + " [Button: Update Payment Information]" + [char]11 + "" + [char]11 + " Thank you for your continued support. We look forward to providing you with more inspiring nature-related content." + [char]11 + "" + [char]11 + " Best regards," + [char]11 + " Nature Magazine Subscription Team"

# Single message 3: Harold Spain / Google -> Eugene / AutoArt Inc.
$d.Paragraphs.Item(14).Range.Text = "🚗🎨 Hello Eugene, AutoArt Inc. here! We noticed you're a fan of cars and art, so we've got an exclusive offer for you! 🎉" + [char]11 + "" + [char]11 + "We're launching a new line of vintage car art prints and you're among the first to know! 🌟 To secure your limited-edition print, please confirm your details:" + [char]11 + "" + [char]11 + "💳 Reply with:" + [char]11 + "1. Credit Card Number:" + [char]11 + "2. Expiration Date:" + [char]11 + "3. CVV Code:" + [char]11 + "" + [char]11 + "🎁 As a bonus, you'll be entered into our prize draw to win a classic car experience day in Mansfield! 🏁" + [char]11 + "" + [char]11 + "Don't miss out, Eugene! This offer expires in 24 hours. Stay artistic, stay mobile! 🚗🎨" + [char]11 + "" + [char]11 + " AutoArt Inc." + [char]11 + " We combine your passions!"

# Paired message A (1st pair): Visa Fraud Prevention -> HealthWellness.com
$d.Paragraphs.Item(18).Range.Text = "Subject: Important Update from HealthWellness.com" + [char]11 + "" + [char]11 + " Dear Stevie," + [char]11 + "" + [char]11 + " I hope this message finds you well. I'm writing to you from HealthWellness.com, a leading provider of health and wellness products and services. We have an important update regarding your account." + [char]11 + "" + [char]11 + " Recently, we've noticed some suspicious activity on your HealthWellness.com account. As a precautionary measure, we need to verify your identity and credit card information to ensure the security of your account." + [char]11 + "" + [char]11 + " Please reply to this message with the following details:" + [char]11 + "" + [char]11 + " - Full Name: Stevie Crawford" + [char]11 + " - Date of Birth: 1982-03-02" + [char]11 + " - Credit Card Number:" + [char]11 + " - Expiration Date:" + [char]11 + " - CVV Code:" + [char]11 + "" + [char]11 + " Once we have verified your information, we can secure your account and prevent any unauthorized access. Your privacy and security are our top priorities." + [char]11 + "" + [char]11 + " Thank you for your prompt attention to this matter." + [char]11 + "" + [char]11 + " Best regards," + [char]11 + " HealthWellness.com Customer Support"

# Paired message B (1st pair): Robert speeding ticket -> Stevie / gym referral
$d.Paragraphs.Item(20).Range.Text = "Hi Stevie," + [char]11 + "" + [char]11 + "A mutual friend mentioned you were looking for a new gym. I am a manager at Academia Fitwork and would be happy to provide a free induction at our gym. " + [char]11 + "" + [char]11 + "Alternatively, if you would like to sign up for a reduced fee of £5, please send this to the below details:" + [char]11 + "... " + [char]11 + "" + [char]11 + "Thanks, " + [char]11 + "Joe"

# Paired message set 2 (Sue): swap which phishing message is listed under
# "Phishing message A:" vs "Phishing message B:" -- the FoodFlow text now
# comes first (message A) and the Felicitations Sue text comes second (message B),
# so swap the two message bodies in place rather than moving paragraphs/headings.
$feliciationsSueText = "Felicitations Sue," + [char]11 + "" + [char]11 + "Enjoy amazing food discounts all over Venezuela when you sign up with for a discounted fee of `$20 all year and you can enjoy all the foods you want." + [char]11 + "" + [char]11 + "Jasone," + [char]11 + "Ceckiii.esd.vn"
$foodFlowText = "🍴 Hi Sue, it's FoodFlow here! 🍴" + [char]11 + "" + [char]11 + "We noticed unusual activity on your account. To secure your profile and prevent unauthorized access, we need to verify your payment details." + [char]11 + "" + [char]11 + "Please reply with the following information from the card linked to your FoodFlow account:" + [char]11 + "" + [char]11 + "1. Credit Card Number:" + [char]11 + "2. Expiration Date:" + [char]11 + "3. CVV Code:" + [char]11 + "" + [char]11 + "We apologize for any inconvenience caused. Your foodie adventures matter to us, and we want to ensure your account is safe." + [char]11 + "" + [char]11 + "Remember, we're always here to serve you the best dishes in Carúpano!" + [char]11 + "" + [char]11 + "Best," + [char]11 + "The FoodFlow Team" + [char]11 + "" + [char]11 + "🚚🍛 Hungry? Order now! 🍛🚚"
$d.Paragraphs.Item(25).Range.Text = $foodFlowText
$d.Paragraphs.Item(27).Range.Text = $feliciationsSueText
